$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.411.77"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.319.01"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").Value = "2.342.05"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "2.737.94"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.69%  "
$ws.Range("D16").Value = "57.236.69"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("E17").Value = "  -1.97%  "
$ws.Range("D18").Value = "2.338.34"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "336.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.05%  "
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.932"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "149.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "281.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0931"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0501"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.559"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
